$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix B64: change from text "4" to numeric 4
$ws.Cells.Item(64, 2).Value = 4

# Add new row 65
$ws.Cells.Item(65, 1).Value = "Ying Tang"

# B65 must stay text "3" (like source data export), not get auto-converted to a number
$ws.Cells.Item(65, 2).NumberFormat = "@"
$ws.Cells.Item(65, 2).Value = "3"
$ws.Cells.Item(65, 2).Style = "Normal"

$ws.Cells.Item(65, 3).Value = "无"
$ws.Cells.Item(65, 4).Value = "FBK"
$ws.Cells.Item(65, 5).Value = "MET"
$ws.Cells.Item(65, 6).Value = "af403c3c-ff8a-4a7e-afb2-ad566d1a3380"
$ws.Cells.Item(65, 7).Value = "ByQZjx-0-_annotated.xlsx"
$ws.Cells.Item(65, 8).Value = "- the use of the ReLU activation, unlike in recurrent highway network where only the tanh activation is used"
